$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H5").Value = 3.3
$ws.Range("I5").Value = 2.42
$ws.Range("J5").Value = 3.2
$ws.Range("L5").Value = 3.05
$ws.Range("N5").Value = 8
$ws.Range("P5").Value = 3.7
$ws.Range("Q5").Value = 1.72
$ws.Range("R5").Value = 2.05
$ws.Range("T5").Value = 2.87
$ws.Range("V5").Value = 2.27
$ws.Range("W5").Value = 11
$ws.Range("X5").Value = 16
$ws.Range("AA5").Value = 20
$ws.Range("AB5").Value = 24
$ws.Range("AC5").Value = 8
$ws.Range("AD5").Value = 6.6
$ws.Range("AG5").Value = 250
$ws.Range("AH5").Value = 9.5
$ws.Range("AI5").Value = 13
$ws.Range("AK5").Value = 26
$ws.Range("AM5").Value = 24
$ws.Range("AN5").Value = 4.8
$ws.Range("AO5").Value = 14.5
$ws.Range("AP5").Value = 19
$ws.Range("AQ5").Value = 60
$ws.Range("AR5").Value = 80
$ws.Range("AS5").Value = 200
$ws.Range("AT5").Value = 2.87
$ws.Range("AW5").Value = 4.5
$ws.Range("AY5").Value = 19.5
$ws.Range("AZ5").Value = 55
$ws.Range("BA5").Value = 80
$ws.Range("G8").Value = 1.42
$ws.Range("H8").Value = 5
$ws.Range("I8").Value = 6
$ws.Range("J8").Value = 1.91
$ws.Range("K8").Value = 2.4
$ws.Range("L8").Value = 6
$ws.Range("M8").Value = 1.03
$ws.Range("N8").Value = 10.5
$ws.Range("U8").Value = 1.83
$ws.Range("V8").Value = 1.83
$ws.Range("X8").Value = 7.5
$ws.Range("Z8").Value = 9.5
$ws.Range("AB8").Value = 23
$ws.Range("AG8").Value = 700
$ws.Range("AH8").Value = 17
$ws.Range("AJ8").Value = 19
$ws.Range("AO8").Value = 7
$ws.Range("AQ8").Value = 19
$ws.Range("AU8").Value = 8.5
$ws.Range("AW8").Value = 8
$ws.Range("AY8").Value = 34
$ws.Range("G9").Value = 1.73
$ws.Range("H9").Value = 3.9
$ws.Range("I9").Value = 4.2
$ws.Range("J9").Value = 2.3
$ws.Range("L9").Value = 4.5
$ws.Range("U9").Value = 1.8
$ws.Range("V9").Value = 1.91
$ws.Range("X9").Value = 8.5
$ws.Range("Y9").Value = 8.5
$ws.Range("Z9").Value = 13
$ws.Range("AA9").Value = 13
$ws.Range("AH9").Value = 13
$ws.Range("AI9").Value = 23
$ws.Range("AJ9").Value = 15
$ws.Range("AL9").Value = 34
$ws.Range("AM9").Value = 41
$ws.Range("AN9").Value = 3.75
$ws.Range("AO9").Value = 9
$ws.Range("AR9").Value = 41
$ws.Range("AW9").Value = 6.5
$ws.Range("AX9").Value = 23
